# TarjetaCreditoPropia.xlsx -- data refresh for automation test fixture
# (commit: "Refactor general de clases, por ejecucion de automatizacion")
#
# Updates the sample rows for the "Datos" sheet: card numbers / card types
# swapped from the old test card (*0702 / *7806, Personal American Express)
# to the new Mastercard test cards (*3585 Empresarial Mastercard,
# *6682 Personal Mastercard), the account id changes from 22493944 to
# 22452521, and rows 5/6 pick up Ahorros account + reference numbers that
# were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 -----------------------------------------------------------
# New account number, entered as a genuine number (General) even though
# the column is displayed with a text format.
$b3 = $ws.Cells.Item(3, 2)
$b3.NumberFormat = "General"
$b3.Value = 22452521
$b3.NumberFormat = "@"

# --- Row 4 -------------------------------------------------------------
$ws.Cells.Item(4, 15).Value = "Empresarial Mastercard"   # O4 tipoTarjeta
$ws.Cells.Item(4, 16).Value = "*3585"                     # P4 numeroTarjeta

# --- Row 5 ---------------------------------------------------------------
$b5 = $ws.Cells.Item(5, 2)
$b5.NumberFormat = "General"
$b5.Value = 22452521
$b5.NumberFormat = "@"

$ws.Cells.Item(5, 15).Value = "Personal Mastercard"      # O5 tipoTarjeta
$ws.Cells.Item(5, 16).Value = "*6682"                      # P5 numeroTarjeta
$ws.Cells.Item(5, 18).Value = "500000"                      # R5 valorPago
$ws.Cells.Item(5, 20).Value = "Ahorros"                      # T5 tipoCuenta
$ws.Cells.Item(5, 21).Value = "406-725210-11"                  # U5 numeroCuenta

# --- Row 6 -----------------------------------------------------------------
$ws.Cells.Item(6, 2).Value = "22452521"                        # B6 (kept as text)

$ws.Cells.Item(6, 15).Value = "Personal Mastercard"              # O6 tipoTarjeta
$ws.Cells.Item(6, 16).Value = "*6682"                              # P6 numeroTarjeta
$ws.Cells.Item(6, 20).Value = "Ahorros"                              # T6 tipoCuenta
$ws.Cells.Item(6, 21).Value = "406-725210-13"                          # U6 numeroCuenta

# --- Selection -------------------------------------------------------------
$ws.Range("U6").Select()
